# Update the "Förändrad" (Changed) date column (C) for all data rows.
# All values in column C (rows 2-119) move forward by one day
# (date serial 45181 -> 45182).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value = 45182
    }
}
